$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unnecessary "nome" (name) column (column B) from the
# "values and proportionality" sheet. Deleting the column shifts all
# subsequent columns (C:M) one position to the left (now B:L).
$ws.Columns.Item(2).Delete()

# Restore the selection to match the saved workbook state.
$ws.Range("K12").Select()
